$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.286.98"
$ws.Range("E2").Value = "'  +0.14%  "
$ws.Range("D3").Value = "'2.285.42"
$ws.Range("E3").Value = "'  -0.84%  "
$ws.Range("E4").Value = "'  +0.02%  "
$ws.Range("D5").Value = "'322.98"
$ws.Range("E5").Value = "'  +1.68%  "
$ws.Range("D6").Value = "'103.12"
$ws.Range("E6").Value = "'  -1.79%  "
$ws.Range("D7").Value = "'0.623"
$ws.Range("E7").Value = "'  -1.00%  "
$ws.Range("E8").Value = "'  +0.08%  "
$ws.Range("D9").Value = "'0.608"
$ws.Range("E9").Value = "'  -0.32%  "
$ws.Range("D10").Value = "'39.85"
$ws.Range("E10").Value = "'  -0.08%  "
$ws.Range("E11").Value = "'  -0.38%  "
$ws.Range("D12").Value = "'8.36"
$ws.Range("E12").Value = "'  -0.58%  "
$ws.Range("D13").Value = "'0.106"
$ws.Range("E13").Value = "'  -0.02%  "
$ws.Range("D14").Value = "'0.972"
$ws.Range("E14").Value = "'  -1.02%  "
$ws.Range("E15").Value = "'  -2.12%  "
$ws.Range("D16").Value = "'2.631.89"
$ws.Range("D17").Value = "'2.285.56"
$ws.Range("E17").Value = "'  -0.95%  "
$ws.Range("D18").Value = "'42.265.29"
$ws.Range("E18").Value = "'  +0.59%  "
$ws.Range("D19").Value = "'7.35"
$ws.Range("E19").Value = "'  -5.54%  "
$ws.Range("E20").Value = "'  -0.71%  "
$ws.Range("D21").Value = "'13.25"
$ws.Range("E21").Value = "'  +31.91%  "
$ws.Range("D22").Value = "'3.62"
$ws.Range("E22").Value = "'  +1.42%  "
$ws.Range("D23").Value = "'73.10"
$ws.Range("E23").Value = "'  -0.78%  "
$ws.Range("D24").Value = "'268.10"
$ws.Range("E24").Value = "'  -7.25%  "
$ws.Range("E25").Value = "'  -2.82%  "
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = "'  -0.09%  "
$ws.Range("D27").Value = "'10.89"
$ws.Range("E27").Value = "'  -0.83%  "
$ws.Range("D28").Value = "'2.30"
$ws.Range("E28").Value = "'  -1.22%  "
$ws.Range("D29").Value = "'22.48"
$ws.Range("E29").Value = "'  -4.22%  "
$ws.Range("D30").Value = "'37.89"
$ws.Range("E30").Value = "'  +6.80%  "
$ws.Range("D31").Value = "'164.05"
$ws.Range("E31").Value = "'  -0.30%  "
$ws.Range("E32").Value = "'  +4.20%  "
$ws.Range("E33").Value = "'  -1.14%  "
$ws.Range("E34").Value = "'  +0.59%  "
$ws.Range("E35").Value = "'  -2.15%  "
$ws.Range("E36").Value = "'  -14.14%  "
$ws.Range("E37").Value = "'  -1.18%  "
$ws.Range("E38").Value = "'  +0.12%  "
$ws.Range("D39").Value = "'3.72"
$ws.Range("E39").Value = "'  +2.52%  "
$ws.Range("E40").Value = "'  -7.36%  "
$ws.Range("E41").Value = "'  +1.68%  "
$ws.Range("D42").Value = "'69.47"
$ws.Range("E42").Value = "'  -2.25%  "
$ws.Range("E43").Value = "'  -0.07%  "
$ws.Range("D44").Value = "'0.225"
$ws.Range("E44").Value = "'  -1.26%  "
$ws.Range("D45").Value = "'91.29"
$ws.Range("D46").Value = "'12.28"
$ws.Range("E46").Value = "'  +1.24%  "
$ws.Range("D47").Value = "'80.41"
$ws.Range("E47").Value = "'  +2.96%  "
$ws.Range("D48").Value = "'112.59"
$ws.Range("E48").Value = "'  -4.05%  "
$ws.Range("D49").Value = "'8.91"
$ws.Range("D51").Value = "'1.596.05"
$ws.Range("E51").Value = "'  +2.82%  "
